# Add a new LeetCode SQL pattern entry ("1321. Restaurant Growth") as a new
# row at the bottom of the Table2 table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question   = "1321. Restaurant Growth"
$difficulty = "Medium"
$pattern    = "Subqueries"
$notes      = "Use window function for rows 6 preceding."
$link       = "https://leetcode.com/problems/restaurant-growth/solutions/1494257/oracle-window-functions-rows-6-preceding/?envType=study-plan-v2&envId=top-sql-50 "

# Grow the table by one row - this also extends the table/autofilter ref
# and the sheet dimension.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

$row = 35

$ws.Range("A$row").Value = $question
$ws.Range("B$row").Value = $difficulty
$ws.Range("C$row").Value = $pattern
$ws.Range("D$row").Value = $notes
$ws.Range("E$row").Value = $link

# Copy the formatting from the row above so the new cells pick up the same
# direct styles (amber "Medium" fill on column B, hyperlink style on E).
$ws.Range("B34").Copy()
$ws.Range("B$row").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("E34").Copy()
$ws.Range("E$row").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Register the hyperlink for the new Link cell.
$ws.Hyperlinks.Add($ws.Range("E$row"), $link)

# Re-apply the hyperlink cell formatting, since adding the hyperlink can
# reset the direct formatting on the cell.
$ws.Range("E34").Copy()
$ws.Range("E$row").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Update the view to match where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("E39").Select()
